# Applies the cryptos.xlsx price/volume refresh described in the commit diff.
# Row 34/35 also swap identity (OKB <-> InjectiveProtocol) with their own new figures.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = '51.362.18'

# Row 3
$ws.Range('D3').Value = '2.974.34'
$ws.Range('E3').Value = '  +1.79%  '

# Row 4
$ws.Range('E4').Value = '  +0.02%  '

# Row 5
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '382.44'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  +1.46%  '

# Row 6
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '103.01'
$ws.Range('D6').Style = "Normal"
$ws.Range('E6').Value = '  +0.22%  '

# Row 7
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '0.541'
$ws.Range('D7').Style = "Normal"
$ws.Range('E7').Value = '  -0.35%  '

# Row 8
$ws.Range('E8').Value = '  +0.02%  '

# Row 9
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.592'
$ws.Range('D9').Style = "Normal"
$ws.Range('E9').Value = '  +1.26%  '

# Row 10
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '36.66'
$ws.Range('D10').Style = "Normal"
$ws.Range('E10').Value = '  -0.84%  '

# Row 11
$ws.Range('E11').Value = '  -0.26%  '

# Row 12
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '0.0842'
$ws.Range('D12').Style = "Normal"
$ws.Range('E12').Value = '  +0.91%  '

# Row 13
$ws.Range('D13').Value = '3.441.46'
$ws.Range('E13').Value = '  +1.94%  '

# Row 14
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '18.16'
$ws.Range('D14').Style = "Normal"
$ws.Range('E14').Value = '  -0.82%  '

# Row 15
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '7.49'
$ws.Range('D15').Style = "Normal"
$ws.Range('E15').Value = '  +2.19%  '

# Row 16
$ws.Range('D16').Value = '2.967.29'
$ws.Range('E16').Value = '  +1.86%  '

# Row 17
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '0.994'
$ws.Range('D17').Style = "Normal"
$ws.Range('E17').Value = '  +7.37%  '

# Row 18
$ws.Range('D18').Value = '51.257.45'
$ws.Range('E18').Value = '  -0.08%  '

# Row 19
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '3.24'
$ws.Range('D19').Style = "Normal"
$ws.Range('E19').Value = '  -4.92%  '

# Row 20
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '7.36'
$ws.Range('D20').Style = "Normal"
$ws.Range('E20').Value = '  +0.25%  '

# Row 21
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '12.72'
$ws.Range('D21').Style = "Normal"
$ws.Range('E21').Value = '  -1.60%  '

# Row 22
$ws.Range('D22').Value = '0.0₃0957'
$ws.Range('E22').Value = '  +1.22%  '

# Row 23
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '68.79'
$ws.Range('D23').Style = "Normal"
$ws.Range('E23').Value = '  +0.71%  '

# Row 24
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '262.41'
$ws.Range('D24').Style = "Normal"
$ws.Range('E24').Value = '  +0.25%  '

# Row 25
$ws.Range('E25').Value = '  +4.88%  '

# Row 26
$ws.Range('E26').Value = '  +12.00%  '

# Row 27
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '7.55'
$ws.Range('D27').Style = "Normal"
$ws.Range('E27').Value = '  +10.16%  '

# Row 28
$ws.Range('E28').Value = '  +12.25%  '

# Row 29
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '4.11'
$ws.Range('D29').Style = "Normal"
$ws.Range('E29').Value = '  -0.24%  '

# Row 30
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '0.166'
$ws.Range('D30').Style = "Normal"
$ws.Range('E30').Value = '  -1.90%  '

# Row 31
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '1.00'
$ws.Range('D31').Style = "Normal"
$ws.Range('E31').Value = '  -0.01%  '

# Row 32
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '25.81'
$ws.Range('D32').Style = "Normal"
$ws.Range('E32').Value = '  +0.60%  '

# Row 33
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '9.86'
$ws.Range('D33').Style = "Normal"
$ws.Range('E33').Value = '  +0.72%  '

# Row 34
$ws.Range('B34').Value = 'InjectiveProtocol'
$ws.Range('C34').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '34.21'
$ws.Range('D34').Style = "Normal"
$ws.Range('E34').Value = '  +0.65%  '

# Row 35
$ws.Range('B35').Value = 'OKB'
$ws.Range('C35').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '50.94'
$ws.Range('D35').Style = "Normal"
$ws.Range('E35').Value = '  -0.50%  '

# Row 36
$ws.Range('E36').Value = '  -2.55%  '

# Row 37
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '0.0446'
$ws.Range('D37').Style = "Normal"
$ws.Range('E37').Value = '  +6.00%  '

# Row 38
$ws.Range('E38').Value = '  -0.02%  '

# Row 39
$ws.Range('E39').Value = '  +0.07%  '

# Row 40
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '17.08'
$ws.Range('D40').Style = "Normal"
$ws.Range('E40').Value = '  +0.98%  '

# Row 41
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '2.56'
$ws.Range('D41').Style = "Normal"
$ws.Range('E41').Value = '  +1.24%  '

# Row 42
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '0.116'
$ws.Range('D42').Style = "Normal"
$ws.Range('E42').Value = '  +1.80%  '

# Row 43
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '1.80'
$ws.Range('D43').Style = "Normal"
$ws.Range('E43').Value = '  -0.89%  '

# Row 44
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '122.37'
$ws.Range('D44').Style = "Normal"
$ws.Range('E44').Value = '  +0.03%  '

# Row 45
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '21.41'
$ws.Range('D45').Style = "Normal"
$ws.Range('E45').Value = '  -0.85%  '

# Row 46
$ws.Range('E46').Value = '  +0.41%  '

# Row 47
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '0.276'
$ws.Range('D47').Style = "Normal"
$ws.Range('E47').Value = '  +2.63%  '

# Row 49
$ws.Range('D49').Value = '2.023.90'
$ws.Range('E49').Value = '  -0.08%  '

# Row 50
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '3.25'
$ws.Range('D50').Style = "Normal"
$ws.Range('E50').Value = '  +3.21%  '

# Row 51
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '0.0337'
$ws.Range('D51').Style = "Normal"
$ws.Range('E51').Value = '  +6.94%  '
